$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at 581, pushing the existing rows 581-640 down to 585-644.
$ws.Rows("581:584").Insert()

# Shared values for the new rows (same market/region/product context as
# every other row in this sheet).
$mercadoId = 9
$mercado   = 'Vega Central Mapocho de Santiago'
$region    = 'Metropolitana'
$codreg    = 13
$tipo      = 'Fruta'
$prodId    = 100108
$producto  = 'Tropicales y subtropicales'
$catId     = 100108006
$categoria = 'Plátano'
$unidad    = '$/caja 20 kilos'
$origen    = 'Ecuador'
$kgUnidad  = 20
$fecha     = 44449

# New row data: Variedad, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm, PrecioKg
$newRows = @(
    @('Barraganete',      'Primera',        450, 19000, 20000, 19444, 972),
    @('Sin especificar',  'Pintón',         680, 18000, 18000, 18000, 900),
    @('Sin especificar',  'Primera Maduro', 550, 19000, 19000, 19000, 950),
    @('Sin especificar',  'Primera Pintón', 600, 20000, 20000, 20000, 1000)
)

$r = 581
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value  = $mercadoId
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = $fecha
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $tipo
    $ws.Cells.Item($r, 7).Value  = $prodId
    $ws.Cells.Item($r, 8).Value  = $producto
    $ws.Cells.Item($r, 9).Value  = $catId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $row[0]
    $ws.Cells.Item($r, 12).Value = $row[1]
    $ws.Cells.Item($r, 13).Value = $row[2]
    $ws.Cells.Item($r, 14).Value = $row[3]
    $ws.Cells.Item($r, 15).Value = $row[4]
    $ws.Cells.Item($r, 16).Value = $row[5]
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $row[6]
    $ws.Cells.Item($r, 20).Value = $kgUnidad
    $r++
}
